# "add more grammar and format changes"
# Applies the grammar/wording fixes from the commit to the resume:
#   1. GPA "3.7" -> "3.70"
#   2. Academic-project bullet: tighten the closing clause about Nvidia
#   3. "offboarded" -> "off-boarded"
#   4. Drop redundant "the " before "development team"
#   5. "design" -> "designing" (and drop the stray trailing period)
#   6. Drop redundant "many " before "company meeting events"

$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $found = $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                                      $true, 1, $false, $replace, 2)
    if (-not $found) {
        Write-Host "WARNING: text not found -> $find"
    }
}

# 1) GPA: "3.7" -> "3.70"
Replace-Text "Cumulative GPA: 3.7" "Cumulative GPA: 3.70"

# 2) Academic project bullet: rewrite ending clause
Replace-Text ("Mao Zedong, Deng Xiaoping, and Xi Jinping and how that has shaped and " + `
              "formed one of the most formidable challenges to Nvidia") `
             "Mao Zedong, Deng Xiaoping, and Xi Jinping and forms a formidable challenge to Nvidia"

# 3) "offboarded" -> "off-boarded"
Replace-Text "Onboarded and offboarded consultants" "Onboarded and off-boarded consultants"

# 4) Remove "the " before "development team"
Replace-Text "Collaborated with the development team and c" "Collaborated with development team and c"

# 5) "design" -> "designing" and drop the trailing period
Replace-Text "design, testing, and implementation of a new CIS system." `
             "designing, testing, and implementation of a new CIS system"

# 6) Remove "many " before "company meeting events"
Replace-Text "Volunteered for many company meeting events" "Volunteered for company meeting events"
